$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay as Text,
# matching the original inline-string / shared-string representation
# (otherwise Excel auto-detects plain numeric-looking strings as numbers).
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D8:D10").NumberFormat = "@"
$ws.Range("D12:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.450.39"
$ws.Range("E2").Value = "  -1.68%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.850.48"
$ws.Range("E3").Value = "  -0.60%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5: BNB
$ws.Range("D5").Value = "243.09"

# Row 6: XRP
$ws.Range("D6").Value = "0.6563"
$ws.Range("E6").Value = "  +2.70%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8: OKB
$ws.Range("D8").Value = "48.13"
$ws.Range("E8").Value = "  +2.72%  "

# Row 9: Cardano
$ws.Range("D9").Value = "0.2998"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "0.07481"
$ws.Range("E10").Value = "  -0.09%  "

# Row 11: Solana
$ws.Range("E11").Value = "  -0.97%  "

# Row 12: TRON
$ws.Range("D12").Value = "0.07632"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.845.01"
$ws.Range("E13").Value = "  -1.01%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "5.016"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15: Polygon
$ws.Range("D15").Value = "0.6849"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "83.72"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "0.000009503"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18: Uniswap
$ws.Range("D18").Value = "6.131"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "29.471.87"
$ws.Range("E19").Value = "  -1.51%  "

# Row 20: WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "2.066.27"
$ws.Range("E20").Value = "  -2.71%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "237.36"
$ws.Range("E21").Value = "  -1.01%  "

# Row 22: Avalanche
$ws.Range("D22").Value = "12.56"
$ws.Range("E22").Value = "  -1.09%  "

# Row 23: Dai
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "7.676"
$ws.Range("E24").Value = "  +3.83%  "

# Row 25: BinanceUSD
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "0.9999"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26: Stellar
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.1426"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27: Monero
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "156.79"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28: Cosmos
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "8.482"
$ws.Range("E28").Value = "  -1.23%  "

# Row 29: EthereumClassic
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "17.80"
$ws.Range("E29").Value = "  -1.09%  "

# Row 30: PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.487"
$ws.Range("E30").Value = "  -1.20%  "

# Row 31: Hedera
$ws.Range("D31").Value = "0.06005"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32: Toncoin
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "1.254"
$ws.Range("E32").Value = "  -2.68%  "

# Row 33: Filecoin
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.135"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.075"
$ws.Range("E34").Value = "  -1.55%  "

# Row 35: LidoDAOToken
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.857"
$ws.Range("E35").Value = "  -1.76%  "

# Row 36: ARBITRUM
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.175"
$ws.Range("E36").Value = "  +0.58%  "

# Row 37: ImmutableX
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.7226"
$ws.Range("E37").Value = "  -1.27%  "

# Row 38: HuobiToken
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.597"
$ws.Range("E38").Value = "  -0.41%  "

# Row 39: MXToken
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  -2.48%  "

# Row 40: VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01779"
$ws.Range("E40").Value = "  -1.61%  "

# Row 41: Maker
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.199.56"
$ws.Range("E41").Value = "  -2.22%  "

# Row 42: FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.235"
$ws.Range("E42").Value = "  -0.72%  "

# Row 43: TrustWalletToken
$ws.Range("D43").Value = "0.9116"
$ws.Range("E43").Value = "  -1.91%  "

# Row 44: PaxDollar
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.23%  "

# Row 45: RocketPoolETH
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.002.81"
$ws.Range("E45").Value = "  -1.56%  "

# Row 46: Quant
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "101.90"
$ws.Range("E46").Value = "  -0.39%  "

# Row 47: Aave
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "66.36"
$ws.Range("E47").Value = "  +0.00%  "

# Row 48: Aptos
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.462"
$ws.Range("E48").Value = "  +10.98%  "

# Row 49: BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.00000000118"
$ws.Range("E49").Value = "  -4.26%  "

# Row 50: TheSandbox
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.4050"
$ws.Range("E50").Value = "  -1.20%  "

# Row 51: EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.100"
$ws.Range("E51").Value = "  -2.79%  "
